$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.606.78"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.739.69"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4912"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D10").Value = "1.745.75"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07047"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.585"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D17").Value = "26.614.77"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007292"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.14%  "
$ws.Range("B19").Value = "BinanceUSD"
$ws.Range("C19").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "1.973.06"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.562"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.713"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.274"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.416"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.025"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08055"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.725"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.613"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6395"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.066"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9043"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.428"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.432"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3918"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.875"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1184"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05394"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.806"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.257"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
